# Adds a new "SiteName" column to the AntennaMetadata sheet (sheet1),
# between the existing AntennaSite (A) and Abbreviation (old B, renamed
# "SiteCode") columns. Values in the new column collapse the per-direction
# AntennaSite names (e.g. "Red Barn (Downstream)") into a shared site
# name (e.g. "Red Barn"). Also corrects A2 from "Below Windy Gap Dam" to
# "Windy Gap Dam".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AntennaMetadata")

# Insert a new column before the old column B (Abbreviation); this shifts
# Abbreviation -> C, UTM_X -> D, UTM_Y -> E.
$ws.Columns.Item(2).Insert()

# Fix a pre-existing typo/inconsistency in the AntennaSite column.
$ws.Cells.Item(2, 1).Value = "Windy Gap Dam"

# New column header + values, entered top-to-bottom so new shared strings
# land in the same order the original author typed them in.
$ws.Cells.Item(1, 2).Value = "SiteName"
$ws.Cells.Item(2, 2).Value = "Windy Gap Dam"
$ws.Cells.Item(3, 2).Value = "Kaibab Park"
$ws.Cells.Item(4, 2).Value = "River Run"
$ws.Cells.Item(5, 2).Value = "Fraser River Canyon"
$ws.Cells.Item(6, 2).Value = "Red Barn"
$ws.Cells.Item(7, 2).Value = "Red Barn"
$ws.Cells.Item(8, 2).Value = "Hitching Post"
$ws.Cells.Item(9, 2).Value = "Hitching Post"
$ws.Cells.Item(10, 2).Value = "Confluence"
$ws.Cells.Item(11, 2).Value = "Confluence"
$ws.Cells.Item(12, 2).Value = "Connectivity Channel Downstream"
$ws.Cells.Item(13, 2).Value = "Connectivity Channel Downstream"
$ws.Cells.Item(14, 2).Value = "Connectivity Channel Side Channel"
$ws.Cells.Item(15, 2).Value = "Connectivity Channel Side Channel"
$ws.Cells.Item(16, 2).Value = "Connectivity Channel Upstream"
$ws.Cells.Item(17, 2).Value = "Connectivity Channel Upstream"

# Rename the old Abbreviation header (now column C) last, so its new
# shared string is appended after all the SiteName strings above.
$ws.Cells.Item(1, 3).Value = "SiteCode"

# Column widths: leave col A (AntennaSite) and col C (old Abbreviation,
# now SiteCode) untouched so their original widths survive unmodified;
# only give the new col B (SiteName) an explicit custom width matching
# col A's (~42.33 characters -- the closest value the character-width/
# pixel-grid rounding allows is 41.42, which lands on the same rendered
# width as col A).
$ws.Columns.Item(2).ColumnWidth = 41.42

# Scroll/selection bookkeeping to match the author's final view state.
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("B6").Select()
